# Applies updated loading_percent values for rows 2-25, columns B,C,D,F,G,H,J
# (columns E, I, K-O remain unchanged at their existing values)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order per data row: B, C, D, F, G, H, J
$cols = @(2, 3, 4, 6, 7, 8, 10)

$data = @(
    ,@(2, 17.21775304652755, 10.62617579077557, 10.14942889989638, 35.01327426648118, 37.32858048598642, 16.10441953963667, 11.36244832300683)
    ,@(3, 16.58010510137931, 10.03089661990075, 10.08169941356973, 34.97795082101561, 37.21939539257713, 16.1599607247745, 11.34755719444668)
    ,@(4, 16.17867488694446, 9.64733113309331, 10.04126599446695, 34.96957913573306, 37.17127819777872, 16.19873762506283, 11.34086337318876)
    ,@(5, 16.01285094746814, 9.486613407096938, 10.02509108268767, 34.96950663612429, 37.15642059558814, 16.21570824189579, 11.33875264070446)
    ,@(6, 15.98518878661125, 9.459664236040126, 10.02242383969566, 34.96969591063399, 37.15424002460661, 16.21859658777038, 11.33843944795332)
    ,@(7, 16.17644723327676, 9.645181308535227, 10.04104661528386, 34.96956465412765, 37.17105860518432, 16.19896177374813, 11.34083240752874)
    ,@(8, 17.0000940692001, 10.42474225432597, 10.12584394486136, 34.99832612826383, 37.2870040289361, 16.12259600168883, 11.35680630280972)
    ,@(9, 18.52671480346259, 11.80566606097868, 10.30072531187987, 35.16064503159132, 37.66443542855507, 16.0102230939921, 11.40749777458103)
    ,@(10, 19.58279255521713, 12.72553429061923, 10.43368719165996, 35.34455726222715, 38.03249598078146, 15.9508354463642, 11.45642944260745)
    ,@(11, 20.04701465424175, 13.12279261295815, 10.49497892194639, 35.44219757500872, 38.21930546113298, 15.9289276326083, 11.48119320955718)
    ,@(12, 20.22033688092132, 13.27013747898108, 10.51829002458448, 35.48116933234974, 38.29279038115479, 15.92137214854584, 11.49092702879948)
    ,@(13, 20.18312066565702, 13.23854201924398, 10.51326529414876, 35.47268746944977, 38.27684279452907, 15.92296632939992, 11.48881489503322)
    ,@(14, 20.06132418966935, 13.13497685447766, 10.49689479594738, 35.44536385799763, 38.22529638120464, 15.92829116370057, 11.48198689909511)
    ,@(15, 19.98639500363154, 13.07113698404804, 10.48688013989404, 35.42888705031172, 38.19407870606314, 15.9316493906202, 11.47785083631217)
    ,@(16, 19.55211630463256, 12.69914235017264, 10.42969660781963, 35.33845652953039, 38.02067403578923, 15.95237047917272, 11.45486115246584)
    ,@(17, 19.28144732807357, 12.46547526697869, 10.39481252529446, 35.28655275278915, 37.91923136440791, 15.96639483327751, 11.44139661140547)
    ,@(18, 19.12424865722585, 12.32908294987729, 10.37482494651657, 35.25801572773785, 37.86271112056648, 15.97494171174586, 11.43388815764789)
    ,@(19, 19.07076791190377, 12.2825619211763, 10.36807113757122, 35.24858003265008, 37.84388925985628, 15.97791786583137, 11.43138656559599)
    ,@(20, 19.31041858374797, 12.49055613912662, 10.39851815102615, 35.29194180781904, 37.92984131317333, 15.96485214869824, 11.44280553659846)
    ,@(21, 20.09716678948719, 13.1654805428372, 10.50170057813413, 35.45333537302039, 38.24036271032733, 15.92670698658996, 11.48398280936997)
    ,@(22, 20.59691094339506, 13.58856970242283, 10.5697195459658, 35.57044970889501, 38.4592767534101, 15.90609586332804, 11.51296934276097)
    ,@(23, 20.33155106134323, 13.36441820959041, 10.53336806849903, 35.50688419559199, 38.3409928827153, 15.91669931575737, 11.49731021926088)
    ,@(24, 19.29732561284339, 12.47922347147658, 10.39684262550825, 35.2895013559296, 37.92503894385061, 15.96554808891858, 11.44216783706887)
    ,@(25, 18.12448618113602, 11.44845333034509, 10.25257270305261, 35.10536944247706, 37.5463043195746, 16.03658046637058, 11.3917213576581)
)

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $col = $cols[$i]
        $val = $entry[$i + 1]
        $ws.Cells.Item($row, $col).Value = $val
    }
}
